# Applies the "Sunny's stack + adils stuff" edit to PseudoCode.docx
# See accompanying unified diff for the target state.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Display Names" + " of people in names.txt" -> single run
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    'Display Names of people in names.txt', $true, $false, $false, $false,
    $false, $true, 1, $false, 'Display Names of people in names.txt', 2)

# ---------------------------------------------------------------------
# 2) "If (" + "$" -> "If ($" (only for the "...primaryChoice)." bullet;
#    there are earlier unrelated "If ($..." occurrences so we locate this
#    one via the unique surrounding phrase first, then narrow the Find to
#    just the two runs that must merge).
# ---------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute(
    'If ($primaryChoice).', $true, $false, $false, $false, $false, $true,
    1, $false, '', 0)
if ($found) {
    $narrow = $d.Range($anchor.Start, $anchor.Start + 5)
    $null = $narrow.Find.Execute('If ($', $true, $false, $false, $false,
        $false, $true, 1, $false, 'If ($', 2)
}

# ---------------------------------------------------------------------
# 3) Remove the old "_GoBack" bookmark (it sat between " has to be
#    Pre-Pended with " and "COUNT...") and merge those two runs back
#    into one, since the split only existed to hold the bookmark.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$mergedTail = ' has to be Pre-Pended with COUNT). First one starts with 1.'
$null = $d.Content.Find.Execute(
    $mergedTail, $true, $false, $false, $false, $false, $true, 1, $false,
    $mergedTail, 2)

# ---------------------------------------------------------------------
# 4) Split "...cannot be part of a pair." into "...cannot be part of a "
#    + "pair." and drop a fresh "_GoBack" bookmark at the join - this is
#    where Word's cursor ended up after the author's last edit.
# ---------------------------------------------------------------------
$tailRange = $d.Content
$found = $tailRange.Find.Execute(
    ', which will contain all combinations which cannot be part of a pair.',
    $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $splitAt = $tailRange.Start + 64   # length of the text up to and incl. "a "
    $bmRange = $d.Range($splitAt, $splitAt)
    $null = $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------
# 5) "Create all pairs for " + "primary - A,A, A-B, A-D, B-A, " -> one run
# ---------------------------------------------------------------------
$dash = [char]0x2013
$createPairsText = 'Create all pairs for primary ' + $dash + ' A,A, A-B, A-D, B-A, '
$null = $d.Content.Find.Execute(
    $createPairsText, $true, $false, $false, $false, $false, $true, 1,
    $false, $createPairsText, 2)

# ---------------------------------------------------------------------
# 6) New bullets: fill the (until now empty) bullet after "... starts
#    with 1." with "Add it to ignoreList", then append three sibling
#    bullets one level back out (ilvl 2 / numId 2).
# ---------------------------------------------------------------------
function Set-BulletText($paragraph, $text, $listLevel) {
    $rng = $paragraph.Range
    $startPos = $rng.Start
    $rng.Text = $text
    if ($listLevel -ne $null) {
        $paragraph.Range.ListFormat.ListLevelNumber = $listLevel
    }
    $fontRange = $d.Range($startPos, $startPos + $text.Length)
    $fontRange.Font.Name = "Arial"
    $fontRange.Font.Size = 10
}

$anchor2 = $d.Content
$null = $anchor2.Find.Execute('First one starts with 1.', $true, $false,
    $false, $false, $false, $true, 1, $false, '', 0)
$anchorPara = $anchor2.Paragraphs(1)
$emptyPara = $anchorPara.Next()

Set-BulletText $emptyPara 'Add it to ignoreList' $null

$emptyPara.Range.InsertParagraphAfter()
$p2 = $emptyPara.Next()
Set-BulletText $p2 'Create a MxM matrix of Names.' 3

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
Set-BulletText $p3 'Remove elements in IgnoreList from the Names Matrix (Set it to ZERO)' 3

$p3.Range.InsertParagraphAfter()
$p4 = $p3.Next()
Set-BulletText $p4 'Enumerate Matrix' 3
